# Updates the "loading_percent" results sheet (res_line) for
# Case_2_213 ("case with 380 kV done") -- refreshed per-line loading
# percentages for rows 2-25 (columns B,C,E,F,G,H,I,K,L,M,O).
# Columns A, D, J, N are unchanged (index / always-zero lines).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 10.33322002669045
$ws.Cells.Item(2, 3).Value = 10.31536733583054
$ws.Cells.Item(2, 5).Value = 13.17156600385459
$ws.Cells.Item(2, 6).Value = 16.86991607391233
$ws.Cells.Item(2, 7).Value = 29.14563638212093
$ws.Cells.Item(2, 8).Value = 14.60368440303813
$ws.Cells.Item(2, 9).Value = 23.28980196499089
$ws.Cells.Item(2, 11).Value = 8.227631653585041
$ws.Cells.Item(2, 12).Value = 10.08628627823765
$ws.Cells.Item(2, 13).Value = 13.12235221975785
$ws.Cells.Item(2, 15).Value = 22.23766179459682
$ws.Cells.Item(3, 2).Value = 10.03282669529225
$ws.Cells.Item(3, 3).Value = 10.3210992370767
$ws.Cells.Item(3, 5).Value = 13.20981206226183
$ws.Cells.Item(3, 6).Value = 15.89584955866808
$ws.Cells.Item(3, 7).Value = 29.33528034154951
$ws.Cells.Item(3, 8).Value = 14.66023726358485
$ws.Cells.Item(3, 9).Value = 23.40507136631799
$ws.Cells.Item(3, 11).Value = 7.975274944797491
$ws.Cells.Item(3, 12).Value = 10.09390744815338
$ws.Cells.Item(3, 13).Value = 13.07050177487505
$ws.Cells.Item(3, 15).Value = 22.34685618342634
$ws.Cells.Item(4, 2).Value = 9.844452854666368
$ws.Cells.Item(4, 3).Value = 10.32534089339875
$ws.Cells.Item(4, 5).Value = 13.23528247850011
$ws.Cells.Item(4, 6).Value = 15.26997757108491
$ws.Cells.Item(4, 7).Value = 29.45985259441859
$ws.Cells.Item(4, 8).Value = 14.69696806203102
$ws.Cells.Item(4, 9).Value = 23.47980731676896
$ws.Cells.Item(4, 11).Value = 7.815000066004924
$ws.Cells.Item(4, 12).Value = 10.09993520498955
$ws.Cells.Item(4, 13).Value = 13.03998707151919
$ws.Cells.Item(4, 15).Value = 22.41801273578925
$ws.Cells.Item(5, 2).Value = 9.766808819655134
$ws.Cells.Item(5, 3).Value = 10.32725169207047
$ws.Cells.Item(5, 5).Value = 13.24616187092072
$ws.Cells.Item(5, 6).Value = 15.00819731993403
$ws.Cells.Item(5, 7).Value = 29.51265783852939
$ws.Cells.Item(5, 8).Value = 14.71244176331161
$ws.Cells.Item(5, 9).Value = 23.51126036520496
$ws.Cells.Item(5, 11).Value = 7.748412626255415
$ws.Cells.Item(5, 12).Value = 10.10273125017459
$ws.Cells.Item(5, 13).Value = 13.02789249315381
$ws.Cells.Item(5, 15).Value = 22.44804412707721
$ws.Cells.Item(6, 2).Value = 9.753866243306057
$ws.Cells.Item(6, 3).Value = 10.32758000812639
$ws.Cells.Item(6, 5).Value = 13.24799859257219
$ws.Cells.Item(6, 6).Value = 14.96433081551593
$ws.Cells.Item(6, 7).Value = 29.52154928759261
$ws.Cells.Item(6, 8).Value = 14.71504172336685
$ws.Cells.Item(6, 9).Value = 23.5165434209114
$ws.Cells.Item(6, 11).Value = 7.737280906521185
$ws.Cells.Item(6, 12).Value = 10.10321606485037
$ws.Cells.Item(6, 13).Value = 13.02590498669813
$ws.Cells.Item(6, 15).Value = 22.45309332437139
$ws.Cells.Item(7, 2).Value = 9.843409135021416
$ws.Cells.Item(7, 3).Value = 10.32536592412542
$ws.Cells.Item(7, 5).Value = 13.23542717694868
$ws.Cells.Item(7, 6).Value = 15.26647399323137
$ws.Cells.Item(7, 7).Value = 29.46055648463072
$ws.Cells.Item(7, 8).Value = 14.69717469741093
$ws.Cells.Item(7, 9).Value = 23.48022746222465
$ws.Cells.Item(7, 11).Value = 7.814107112633764
$ws.Cells.Item(7, 12).Value = 10.09997153726243
$ws.Cells.Item(7, 13).Value = 13.03982257108938
$ws.Cells.Item(7, 15).Value = 22.41841355976832
$ws.Cells.Item(8, 2).Value = 10.23052292034961
$ws.Cells.Item(8, 3).Value = 10.31719419978901
$ws.Cells.Item(8, 5).Value = 13.18434098777264
$ws.Cells.Item(8, 6).Value = 16.53996406344768
$ws.Cells.Item(8, 7).Value = 29.20933618734731
$ws.Cells.Item(8, 8).Value = 14.62276779339808
$ws.Cells.Item(8, 9).Value = 23.32872599462207
$ws.Cells.Item(8, 11).Value = 8.141761998982458
$ws.Cells.Item(8, 12).Value = 10.08863461246852
$ws.Cells.Item(8, 13).Value = 13.10420480458766
$ws.Cells.Item(8, 15).Value = 22.27445933183361
$ws.Cells.Item(9, 2).Value = 10.95409679206852
$ws.Cells.Item(9, 3).Value = 10.30687069590614
$ws.Cells.Item(9, 5).Value = 13.09991498607624
$ws.Cells.Item(9, 6).Value = 19.00274580682531
$ws.Cells.Item(9, 7).Value = 28.781354362782
$ws.Cells.Item(9, 8).Value = 14.49274058274427
$ws.Cells.Item(9, 9).Value = 23.06297324753785
$ws.Cells.Item(9, 11).Value = 8.739457238025599
$ws.Cells.Item(9, 12).Value = 10.0770715140107
$ws.Cells.Item(9, 13).Value = 13.24057810088145
$ws.Cells.Item(9, 15).Value = 22.02474438669891
$ws.Cells.Item(10, 2).Value = 11.45874993441797
$ws.Cells.Item(10, 3).Value = 10.30272226187399
$ws.Cells.Item(10, 5).Value = 13.04747227879457
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 28.5065466784079
$ws.Cells.Item(10, 8).Value = 14.40683567768794
$ws.Cells.Item(10, 9).Value = 22.8867183644111
$ws.Cells.Item(10, 11).Value = 9.148215421532852
$ws.Cells.Item(10, 12).Value = 10.07503722478731
$ws.Cells.Item(10, 13).Value = 13.34643188105195
$ws.Cells.Item(10, 15).Value = 21.86108644826114
$ws.Cells.Item(11, 2).Value = 11.68151796710968
$ws.Cells.Item(11, 3).Value = 10.30157262954614
$ws.Cells.Item(11, 5).Value = 13.02569186532505
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 28.39018331669142
$ws.Cells.Item(11, 8).Value = 14.3698334951215
$ws.Cells.Item(11, 9).Value = 22.81063625668751
$ws.Cells.Item(11, 11).Value = 9.327041524159299
$ws.Cells.Item(11, 12).Value = 10.07550423834301
$ws.Cells.Item(11, 13).Value = 13.39570257448745
$ws.Cells.Item(11, 15).Value = 21.79092474893183
$ws.Cells.Item(12, 2).Value = 11.76482676265953
$ws.Cells.Item(12, 3).Value = 10.30124260104744
$ws.Cells.Item(12, 5).Value = 13.01774239103026
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 28.347367576779
$ws.Cells.Item(12, 8).Value = 14.35611945392551
$ws.Cells.Item(12, 9).Value = 22.78241345208311
$ws.Cells.Item(12, 11).Value = 9.393696344168548
$ws.Cells.Item(12, 12).Value = 10.07588028192981
$ws.Cells.Item(12, 13).Value = 13.41451072486552
$ws.Cells.Item(12, 15).Value = 21.7649722166802
$ws.Cells.Item(13, 2).Value = 11.74693240858431
$ws.Cells.Item(13, 3).Value = 10.30130900672666
$ws.Cells.Item(13, 5).Value = 13.019441189677
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 28.35653310268737
$ws.Cells.Item(13, 8).Value = 14.35905978063649
$ws.Cells.Item(13, 9).Value = 22.78846561437482
$ws.Cells.Item(13, 11).Value = 9.379388873443776
$ws.Cells.Item(13, 12).Value = 10.07579045221016
$ws.Cells.Item(13, 13).Value = 13.41045352913693
$ws.Cells.Item(13, 15).Value = 21.7705341583188
$ws.Cells.Item(14, 2).Value = 11.68839322667463
$ws.Cells.Item(14, 3).Value = 10.3015433715746
$ws.Cells.Item(14, 5).Value = 13.02503188040025
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 28.3866357884436
$ws.Cells.Item(14, 8).Value = 14.36869926663999
$ws.Cells.Item(14, 9).Value = 22.80830257583709
$ws.Cells.Item(14, 11).Value = 9.3325467620691
$ws.Cells.Item(14, 12).Value = 10.07553118968978
$ws.Cells.Item(14, 13).Value = 13.39724697739072
$ws.Cells.Item(14, 15).Value = 21.78877727085408
$ws.Cells.Item(15, 2).Value = 11.6523977533327
$ws.Cells.Item(15, 3).Value = 10.30170061938392
$ws.Cells.Item(15, 5).Value = 13.02849517956272
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 28.40523729925003
$ws.Cells.Item(15, 8).Value = 14.37464250274764
$ws.Cells.Item(15, 9).Value = 22.82052980288337
$ws.Cells.Item(15, 11).Value = 9.303715138434882
$ws.Cells.Item(15, 12).Value = 10.07539829308578
$ws.Cells.Item(15, 13).Value = 13.38917686287722
$ws.Cells.Item(15, 15).Value = 21.80003194085802
$ws.Cells.Item(16, 2).Value = 11.44404845895978
$ws.Cells.Item(16, 3).Value = 10.30281217060945
$ws.Cells.Item(16, 5).Value = 13.04893743532983
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 28.51432573914338
$ws.Cells.Item(16, 8).Value = 14.4092955793627
$ws.Cells.Item(16, 9).Value = 22.89177283962167
$ws.Cells.Item(16, 11).Value = 9.136382013262359
$ws.Cells.Item(16, 12).Value = 10.07503463692604
$ws.Cells.Item(16, 13).Value = 13.34323353593859
$ws.Cells.Item(16, 15).Value = 21.86575792113124
$ws.Cells.Item(17, 2).Value = 11.31443839289
$ws.Cells.Item(17, 3).Value = 10.30368244234558
$ws.Cells.Item(17, 5).Value = 13.06200965612716
$ws.Cells.Item(17, 6).Value = 20.20408069597325
$ws.Cells.Item(17, 7).Value = 28.58346641626823
$ws.Cells.Item(17, 8).Value = 14.4310854289417
$ws.Cells.Item(17, 9).Value = 22.93652655196216
$ws.Cells.Item(17, 11).Value = 9.031875819121913
$ws.Cells.Item(17, 12).Value = 10.07516746810694
$ws.Cells.Item(17, 13).Value = 13.3153277249505
$ws.Cells.Item(17, 15).Value = 21.90717648404797
$ws.Cells.Item(18, 2).Value = 11.23925369433529
$ws.Cells.Item(18, 3).Value = 10.30425246276607
$ws.Cells.Item(18, 5).Value = 13.06972387314139
$ws.Cells.Item(18, 6).Value = 19.95656407809801
$ws.Cells.Item(18, 7).Value = 28.62404799261193
$ws.Cells.Item(18, 8).Value = 14.44381385444852
$ws.Cells.Item(18, 9).Value = 22.96265342064629
$ws.Cells.Item(18, 11).Value = 8.971098470981262
$ws.Cells.Item(18, 12).Value = 10.07537499715041
$ws.Cells.Item(18, 13).Value = 13.29938283170564
$ws.Cells.Item(18, 15).Value = 21.93140288461846
$ws.Cells.Item(19, 2).Value = 11.21369041344482
$ws.Cells.Item(19, 3).Value = 10.30445741461612
$ws.Cells.Item(19, 5).Value = 13.07236934525925
$ws.Cells.Item(19, 6).Value = 19.87204792380568
$ws.Cells.Item(19, 7).Value = 28.63792781991889
$ws.Cells.Item(19, 8).Value = 14.44815707618283
$ws.Cells.Item(19, 9).Value = 22.97156581773344
$ws.Cells.Item(19, 11).Value = 8.95040680978425
$ws.Cells.Item(19, 12).Value = 10.07546781448466
$ws.Cells.Item(19, 13).Value = 13.29400264360562
$ws.Cells.Item(19, 15).Value = 21.93967484577119
$ws.Cells.Item(20, 2).Value = 11.32830200698474
$ws.Cells.Item(20, 3).Value = 10.30358261659101
$ws.Cells.Item(20, 5).Value = 13.06059787220464
$ws.Cells.Item(20, 6).Value = 20.24955283636154
$ws.Cells.Item(20, 7).Value = 28.57602202407085
$ws.Cells.Item(20, 8).Value = 14.42874563589955
$ws.Cells.Item(20, 9).Value = 22.93172253233297
$ws.Cells.Item(20, 11).Value = 9.0430701289271
$ws.Cells.Item(20, 12).Value = 10.07513976397682
$ws.Cells.Item(20, 13).Value = 13.31828747081446
$ws.Cells.Item(20, 15).Value = 21.90272564492826
$ws.Cells.Item(21, 2).Value = 11.70561658976062
$ws.Cells.Item(21, 3).Value = 10.30147168075836
$ws.Cells.Item(21, 5).Value = 13.023381665304
$ws.Cells.Item(21, 6).Value = 21.46857628470577
$ws.Cells.Item(21, 7).Value = 28.37775997463276
$ws.Cells.Item(21, 8).Value = 14.36585983698719
$ws.Cells.Item(21, 9).Value = 22.80246003596986
$ws.Cells.Item(21, 11).Value = 9.3463345558598
$ws.Cells.Item(21, 12).Value = 10.07560194388082
$ws.Cells.Item(21, 13).Value = 13.40112206246215
$ws.Cells.Item(21, 15).Value = 21.78340210567103
$ws.Cells.Item(22, 2).Value = 11.94607355196682
$ws.Cells.Item(22, 3).Value = 10.30070555041318
$ws.Cells.Item(22, 5).Value = 13.00079729884418
$ws.Cells.Item(22, 6).Value = 22.22866616901552
$ws.Cells.Item(22, 7).Value = 28.25546411804202
$ws.Cells.Item(22, 8).Value = 14.32649630587088
$ws.Cells.Item(22, 9).Value = 22.72140531036894
$ws.Cells.Item(22, 11).Value = 9.53832596543411
$ws.Cells.Item(22, 12).Value = 10.07706453588847
$ws.Cells.Item(22, 13).Value = 13.45613105745606
$ws.Cells.Item(22, 15).Value = 21.70900874320163
$ws.Cells.Item(23, 2).Value = 11.81832043306158
$ws.Cells.Item(23, 3).Value = 10.30105856235466
$ws.Cells.Item(23, 5).Value = 13.01269200178153
$ws.Cells.Item(23, 6).Value = 21.82633154458858
$ws.Cells.Item(23, 7).Value = 28.32006793271338
$ws.Cells.Item(23, 8).Value = 14.34734674819278
$ws.Cells.Item(23, 9).Value = 22.76435270914133
$ws.Cells.Item(23, 11).Value = 9.436436401634923
$ws.Cells.Item(23, 12).Value = 10.07617809573267
$ws.Cells.Item(23, 13).Value = 13.42669539017128
$ws.Cells.Item(23, 15).Value = 21.74838534766764
$ws.Cells.Item(24, 2).Value = 11.32203634896872
$ws.Cells.Item(24, 3).Value = 10.30362753068985
$ws.Cells.Item(24, 5).Value = 13.06123552039328
$ws.Cells.Item(24, 6).Value = 20.22900810905287
$ws.Cells.Item(24, 7).Value = 28.57938504354707
$ws.Cells.Item(24, 8).Value = 14.42980282971146
$ws.Cells.Item(24, 9).Value = 22.9338931917531
$ws.Cells.Item(24, 11).Value = 9.03801134432916
$ws.Cells.Item(24, 12).Value = 10.07515188038362
$ws.Cells.Item(24, 13).Value = 13.31694906202624
$ws.Cells.Item(24, 15).Value = 21.90473657870628
$ws.Cells.Item(25, 2).Value = 10.76273359063799
$ws.Cells.Item(25, 3).Value = 10.30905714351201
$ws.Cells.Item(25, 5).Value = 13.12106981171672
$ws.Cells.Item(25, 6).Value = 18.34778573295695
$ws.Cells.Item(25, 7).Value = 28.89019111542679
$ws.Cells.Item(25, 8).Value = 14.52622185453523
$ws.Cells.Item(25, 9).Value = 23.13152292330444
$ws.Cells.Item(25, 11).Value = 8.582899502895765
$ws.Cells.Item(25, 12).Value = 10.07906161466627
$ws.Cells.Item(25, 13).Value = 13.20265350156913
$ws.Cells.Item(25, 15).Value = 22.08881679889081
